$d = $word.ActiveDocument

# Bold + colored highlight applied to quantitative impact metrics (hex 2C3E50,
# expressed as the decimal BGR-style value the Word COM Font.Color property expects).
$metricColor = 5258796

function Highlight-Metric($rng, $text) {
    $found = $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Bold = 1
        $rng.Font.Color = $metricColor
    }
    return $found
}

function Assert-ParagraphContains($rng, $marker) {
    if ($rng.Text -notlike "*$marker*") {
        throw "Paragraph did not contain expected marker: $marker"
    }
}

# "• Discovered systematic race coding errors ... improving demographic classification
#  accuracy from 23% to 64%" -> bold/color "23%" and "64%"
$rng = $d.Paragraphs.Item(9).Range
Assert-ParagraphContains $rng "23% to 64%"
Highlight-Metric $rng "23%"
Highlight-Metric $rng "64%"

# "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%,
#  reducing polling error margins from ±4.2% to ±2.1%"
# -> bold/color "87%", "71%", "±4.2%" and "±2.1%"
$rng = $d.Paragraphs.Item(11).Range
Assert-ParagraphContains $rng "reducing polling error margins"
Highlight-Metric $rng "87%"
Highlight-Metric $rng "71%"
Highlight-Metric $rng "±4.2%"
Highlight-Metric $rng "±2.1%"

# "• Wrote RFP and analyzed bids from 1,200 vendors for research platform development"
# -> bold/color "1,200"
$rng = $d.Paragraphs.Item(31).Range
Assert-ParagraphContains $rng "1,200 vendors"
Highlight-Metric $rng "1,200"

# "• Created comprehensive meta-analysis framework ... became the $400M Polling
#  Consortium Database at The Analyst Institute, now valued at $1B+"
# -> bold/color "$400M" and "$1B" (the trailing "+" stays plain)
$rng = $d.Paragraphs.Item(46).Range
Assert-ParagraphContains $rng "Polling Consortium Database"
Highlight-Metric $rng "$400M"
Highlight-Metric $rng "$1B"

# "• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M"
# -> bold/color "73.5%" and "$4.7M"
$rng = $d.Paragraphs.Item(63).Range
Assert-ParagraphContains $rng "campaigns and organizations"
Highlight-Metric $rng "73.5%"
Highlight-Metric $rng "$4.7M"

# "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%"
# -> bold/color "87%" and "71%"
$rng = $d.Paragraphs.Item(65).Range
Assert-ParagraphContains $rng "industry standard of 71%"
Highlight-Metric $rng "87%"
Highlight-Metric $rng "71%"
